$d = $word.ActiveDocument

$replacements = @(
    @{old="852÷3="; new="282÷7="},
    @{old="243÷5="; new="697÷8="},
    @{old="107÷8="; new="560÷3="},
    @{old="990÷4="; new="168÷9="},
    @{old="392÷8="; new="181÷8="},
    @{old="775÷6="; new="257÷6="},
    @{old="992÷7="; new="541÷8="},
    @{old="862÷2="; new="177÷7="},
    @{old="497÷5="; new="731÷9="},
    @{old="755÷2="; new="125÷3="},
    @{old="146÷8="; new="268÷2="},
    @{old="361÷7="; new="722÷9="},
    @{old="691÷9="; new="474÷9="},
    @{old="221÷2="; new="532÷4="},
    @{old="648÷9="; new="522÷4="},
    @{old="976÷3="; new="583÷5="},
    @{old="234÷8="; new="656÷2="},
    @{old="364÷3="; new="530÷4="},
    @{old="403÷8="; new="476÷6="},
    @{old="196÷4="; new="829÷7="},
    @{old="141÷3="; new="102÷7="},
    @{old="723÷2="; new="275÷5="},
    @{old="536÷6="; new="829÷4="},
    @{old="665÷6="; new="993÷3="},
    @{old="367÷7="; new="723÷5="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
